# Bump the "Förändrad" (Changed) date column (C) by one day for every data row
# (rows 2-463): serial date 45181 (2023-09-12) -> 45182 (2023-09-13)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($r = 2; $r -le 463; $r++) {
    $ws.Cells.Item($r, 3).Value = 45182
}
